$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LUT): Utilization and Utilization % updated
$ws.Range("B2").Value = 11169.0
$ws.Range("D2").Value = 63.46022415161133

# Row 3 (LUTRAM): Utilization and Utilization % updated
$ws.Range("B3").Value = 655.0
$ws.Range("D3").Value = 10.916666984558105

# Row 4 (FF): Utilization and Utilization % updated
$ws.Range("B4").Value = 13982.0
$ws.Range("D4").Value = 39.72159194946289
